# The workbook currently has two sheets in this order:
#   1) "2022-Q2"  - the big per-fund holdings table (A1:H88)
#   2) "总计"      - the small summary table (A1:D2)
#
# The target state keeps the same two tab positions, but swaps which data
# (and name) lives in each position:
#   1) "总计"      - the small summary table (A1:D2)
#   2) "2022-Q2"  - the big per-fund holdings table (A1:H88)
#
# This is done as a pure content/name swap (no tab reordering, no adding or
# removing sheets - inserting/removing a sheet re-indexes already-captured
# worksheet handles in this host, so we stick to the original two sheets
# throughout) - only the name + data each position shows changes.

$wb = $excel.ActiveWorkbook

$wsBig = $wb.Worksheets.Item(1)   # currently "2022-Q2"   (A1:H88)
$wsSmall = $wb.Worksheets.Item(2) # currently "总计"        (A1:D2)

# 1) Stage the small table (values + formatting) off in unused columns of
#    the big sheet, well clear of its A1:H88 data, so it isn't disturbed by
#    the next steps.
$wsSmall.Range("A1:D2").Copy($wsBig.Range("Z1"))

# 2) The small sheet's real data is now staged elsewhere, so wipe it.
$wsSmall.Cells.Clear()

# 3) Move the big table into what will become the "2022-Q2" sheet.
$wsBig.Range("A1:H88").Copy($wsSmall.Range("A1"))

# 4) Clear the big table out of its old home (but leave the staged small
#    table in the Z1:AC2 corner alone).
$wsBig.Range("A1:H88").Clear()

# 5) Drop the staged small table into place as the new A1:D2, then clear
#    the staging area.
$wsBig.Range("Z1:AC2").Copy($wsBig.Range("A1"))
$wsBig.Range("Z1:AC2").Clear()

# 6) Swap the two sheet names. Rename through a temporary name first so we
#    never try to give a sheet the name another sheet still currently has.
$wsBig.Name = "__tmp_rename__"
$wsSmall.Name = "2022-Q2"
$wsBig.Name = "总计"

# 7) The big "2022-Q2" sheet is the selected/active tab, same as before.
$wsSmall.Select()
